$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Req. ID"
$ws.Range("E3").Value = "Sub-module"
